$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 2851.6667
$ws.Cells.Item(2, 9).Value = 852.5
$ws.Cells.Item(2, 11).Value = 852.5
$ws.Cells.Item(2, 13).Value = -739.5

# ALC row 20
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 1229.6666
$ws.Cells.Item(20, 9).Value = 1229.6666
$ws.Cells.Item(20, 11).Value = 1229.6666
$ws.Cells.Item(20, 13).Value = -999.6666

# ALC row 35
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(35, 8).Value = 1229.6666
$ws.Cells.Item(35, 9).Value = 1229.6666
$ws.Cells.Item(35, 11).Value = 1229.6666
$ws.Cells.Item(35, 13).Value = -850.6666

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3499.5
$ws.Cells.Item(43, 9).Value = 3499.5
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 3499.5
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -3430.5
$ws.Cells.Item(43, 14).ClearContents()

# ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(47, 8).Value = 10000
$ws.Cells.Item(47, 10).Value = 10000
$ws.Cells.Item(47, 12).Value = 10000
$ws.Cells.Item(47, 14).Value = -11944

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 297.5
$ws.Cells.Item(131, 10).Value = 295
$ws.Cells.Item(131, 12).Value = 885
$ws.Cells.Item(131, 14).Value = -10965

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1757.4138
$ws.Cells.Item(132, 9).Value = 1757.4138
$ws.Cells.Item(132, 11).Value = 5272.2414
$ws.Cells.Item(132, 13).Value = -2742.2414

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 2694.25
$ws.Cells.Item(110, 9).Value = 388.5
$ws.Cells.Item(110, 11).Value = 388.5
$ws.Cells.Item(110, 13).Value = 1656.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2789.8
$ws.Cells.Item(122, 9).Value = 1914.5714
$ws.Cells.Item(122, 10).Value = 4832
$ws.Cells.Item(122, 11).Value = 5743.7142
$ws.Cells.Item(122, 12).Value = 14496
$ws.Cells.Item(122, 13).Value = -3293.7142
$ws.Cells.Item(122, 14).Value = -19396

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 920
$ws.Cells.Item(132, 9).Value = 847.5
$ws.Cells.Item(132, 11).Value = 2542.5
$ws.Cells.Item(132, 13).Value = -12.5

# BSM row 141
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(141, 8).Value = 44166.668
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 44166.668
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 44166.668
$ws.Cells.Item(141, 14).Value = -54526.668
$ws.Cells.Item(141, 13).ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2718.4614
$ws.Cells.Item(58, 9).Value = 2528.3333
$ws.Cells.Item(58, 10).Value = 5000
$ws.Cells.Item(58, 11).Value = 2528.3333
$ws.Cells.Item(58, 12).Value = 5000
$ws.Cells.Item(58, 13).Value = -2325.3333
$ws.Cells.Item(58, 14).Value = -5406

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 5338.75
$ws.Cells.Item(99, 9).Value = 3800
$ws.Cells.Item(99, 11).Value = 3800
$ws.Cells.Item(99, 13).Value = -2302

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 5338.75
$ws.Cells.Item(126, 9).Value = 3800
$ws.Cells.Item(126, 11).Value = 11400
$ws.Cells.Item(126, 13).Value = -8930

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2873.75
$ws.Cells.Item(132, 9).Value = 2570
$ws.Cells.Item(132, 11).Value = 7710
$ws.Cells.Item(132, 13).Value = -5180

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2718.4614
$ws.Cells.Item(136, 9).Value = 2528.3333
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 7584.999899999999
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -5034.999899999999
$ws.Cells.Item(136, 14).Value = -20100

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 13).ClearContents()

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 35
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 14).ClearContents()

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 595
$ws.Cells.Item(55, 9).Value = 595
$ws.Cells.Item(55, 11).Value = 1785
$ws.Cells.Item(55, 13).Value = -1608

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 500
$ws.Cells.Item(109, 9).Value = 500
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 1500
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -460
$ws.Cells.Item(109, 14).ClearContents()

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 4550.5
$ws.Cells.Item(132, 9).Value = 2999.5
$ws.Cells.Item(132, 10).Value = 7652.5
$ws.Cells.Item(132, 11).Value = 26995.5
$ws.Cells.Item(132, 12).Value = 68872.5
$ws.Cells.Item(132, 13).Value = -24465.5
$ws.Cells.Item(132, 14).Value = -73932.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2404.3635
$ws.Cells.Item(113, 9).Value = 2082.3333
$ws.Cells.Item(113, 10).Value = 2525.125
$ws.Cells.Item(113, 11).Value = 2082.3333
$ws.Cells.Item(113, 12).Value = 2525.125
$ws.Cells.Item(113, 13).Value = 87.66670000000022
$ws.Cells.Item(113, 14).Value = -6865.125

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1494.3334
$ws.Cells.Item(132, 9).Value = 1494.3334
$ws.Cells.Item(132, 11).Value = 4483.0002
$ws.Cells.Item(132, 13).Value = -1953.0002

# LTW row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 969.5
$ws.Cells.Item(9, 9).Value = 608
$ws.Cells.Item(9, 10).Value = 3500
$ws.Cells.Item(9, 11).Value = 608
$ws.Cells.Item(9, 12).Value = 3500
$ws.Cells.Item(9, 13).Value = -384
$ws.Cells.Item(9, 14).Value = -3948

# LTW row 30
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 290.75
$ws.Cells.Item(30, 9).Value = 290.75
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 290.75
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = -182.75
$ws.Cells.Item(30, 14).ClearContents()

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3200
$ws.Cells.Item(46, 9).Value = 2975.25
$ws.Cells.Item(46, 11).Value = 2975.25
$ws.Cells.Item(46, 13).Value = -2787.25

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 874.5
$ws.Cells.Item(93, 10).Value = 874.5
$ws.Cells.Item(93, 12).Value = 874.5
$ws.Cells.Item(93, 14).Value = -3370.5

# LTW row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(107, 8).Value = 4000
$ws.Cells.Item(107, 9).Value = 4000
$ws.Cells.Item(107, 11).Value = 4000
$ws.Cells.Item(107, 13).Value = -2080

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 7903.4375
$ws.Cells.Item(132, 9).Value = 7834.143
$ws.Cells.Item(132, 10).Value = 8388.5
$ws.Cells.Item(132, 11).Value = 23502.429
$ws.Cells.Item(132, 12).Value = 25165.5
$ws.Cells.Item(132, 13).Value = -20972.429
$ws.Cells.Item(132, 14).Value = -30225.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2846.1
$ws.Cells.Item(132, 9).Value = 2884.6667
$ws.Cells.Item(132, 10).Value = 2499
$ws.Cells.Item(132, 11).Value = 8654.000100000001
$ws.Cells.Item(132, 12).Value = 7497
$ws.Cells.Item(132, 13).Value = -6124.000100000001
$ws.Cells.Item(132, 14).Value = -12557
